$d = $word.ActiveDocument

# 1) Date line
$null = $d.Content.Find.Execute("המאמר היומי של מייק: 26.07.25", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר היומי של מייק: 23.07.25", 2)

# 2) Title
$null = $d.Content.Find.Execute("Building Bridges between Regression, Clustering, and Classification", $true, $false, $false, $false, $false, $true, 1, $false, "Reinforcement Pre-Training", 2)

# 3) Intro paragraph
$null = $d.Content.Find.Execute("מזמן לא סקרתי מאמר שלא מופיעה בו גם מילה LLM וגם diffusion models - תתפלאו אבל יש עדיין כאלו ואני חייב להודות שזה היווה אחת הסיבות לבחירתו. המאמר דן בבעיה די מעניינת היא המרה של בעיות רגרסיה לבעיות סיווג (בתחום למידה עמוקה). ", $true, $false, $false, $false, $false, $true, 1, $false, "חוזר מחופשה עם סקירה מאוד קצרה של רעיון מאוד מסקרן ודי אינטואיטיבי לאימון של מודל שפה. אנו רגילים שבשלב הראשון של אימון מודל שפה, הנקרא אימון מקדים, אנו מאמנים אותו על מה שנקרא  next token prediction או NTP. כלומר בהינתן דאטהסט עצום ולא מתויג אנו ממקסמים את הנראות (likelihood) עבור כל טוקן בדאטהסט בהינתן ההקשר שלו כלומר כל הטוקנים לפניו. המטרה כאן היא למקסם את הנראות המשוערכת של הדאטהסט עם המודל המאומן (ניתן לראות זאת באמצעות שימוש פשוט בחוק בייס). ד״א ניתן לראות די בקלות שבאמצעות אימון מקדים כזה המודל מוסגל לרכוש מיומנויות רבות כלומר ידע במגוון תחומים, פתרון שאלות פשוטות וכדומה. ", 2)

# 4) "most deep models..." paragraph
$null = $d.Content.Find.Execute("מרבית המודלים העמוקים שלנו היום, כמו llms, מודלים ויזואליים ומולטימודליים הם מודלי סיווג במהותם כלומר הפלט שלהם חי במרחב דיסקרטי כלשהו למשל טוענים טקסטואליים או פיקסלים. אז זה נשמע די טבעי לקחת בעיה שהפלט שלה רציף (חד או רב מימדי), להמיר אותה לבעיית סיווג ולבנות (לאמן) מודל סיווג במקום מודל רגרסיה. זה נעשה בד״כ על ידי חלוקה(binning) של מרחב הפלט לכמה תת-מרחבים זרים ואז כל פלט ממופה למספר תת-המרחב שהוא שייך אליו. ככה בעיית רגרסיה הופכת להיות בעיית סיווג. לאחר אימון המודל ניתן להמיר את הערך הדיסקרטי בחזרה למרחב הרציף תוך שימוש חיזוי המודל (לרוב סופטמקס).", $true, $false, $false, $false, $false, $true, 1, $false, "אחרי השלב הראשון באים השלבים של alignment כלומר SFT שזה Supervised Fine Tuning וגם RLHF (עם כל סדר ביניהם). המאמר שסוקרים היום שואל את השאלה הבא: למה לא ניתן לבצע אימון NTP על כל הדאטהסט עם למידה עם חיזוקים או RL. מתברר שזה אפשרי ויש לזה פוטנציאל לשיפור ביצועי המודל.", 2)

# 5) "the paper we review today proposes..." paragraph
$null = $d.Content.Find.Execute("המאמר שנסקור היום מציע גישה כללית לפיתוח מודלי סיווג לבעיות רציפות. המחברים מציעים כמה מודלים שמאומנים בצוותא לפתרון בעיה זו. המודל הראשון, האנקודר, לוקח את הקלט מעביר אותה למרחב הלטנטי ובנוסף מאמנים שכבה שחוזה את התפלגות הקטגוריות עבור הקלט (אחרי ההמרה). ", $true, $false, $false, $false, $false, $true, 1, $false, "איך עושים זאת בפועל? עבור כל טוקן בטקסט אנו מבקשים מהמודל לעשות תהליך ריזונינג קצר כדי לנחש את הטוקן הבא. המודל מתבקש ליצור כמה מסלולי חשיבה כאלו - המסלול שמנחש את המילה בצורה נכונה מקבל תגמול 1 כאשר השאר מקבלים 0. לאחר מכן ניתן להשתמש בתגמולים אלו כדי לאמן את המודל בשיטה האהובה שלכם מעולם למידה עם חיזוקים (PPO, GRPO וכל שיטה אחרת). כלומר מקרה די קלאסי של RLVR או שזה RL עם verifiable rewards.", 2)

# 6) "the second model..." paragraph
$null = $d.Content.Find.Execute("המודל השני לוקח את הפלט ומעביר אותו למרחב החדש של הקטגוריות. הקטגוריה של הפלט יכולה להיות רכה או soft - כלומר להוות התפלגות לא מנוונת(לא וקטור one-hot) מעל כל הקטגוריות. משמעות הדבר שהתפלגות יעד של הקטגוריה עבור פלטים מסוימים, הקרובים לכמה מרכזי קלסטרים, תשקף את זה בצורה הסתברותית. מה שמאומן במודל הזה זה מרכזי הקלסטרים. התפלגות קטגוריות עבור הפלט מחושבת למשל עם פונקציית סופטמקס המשקללת את הסיכוי של הפלט שייך לקלסטר המחושב באמצעות התפלגות גאוסית (למשל). שני המודלים אלו מאומנים יחד כאשר פונקציית לוס הוא מרחק KL בין התפלגויות הקטגוריות שהן מוציאות.", $true, $false, $false, $false, $false, $true, 1, $false, "ההבדל העיקרי בין שיטת אימון מקדים זו ל-pretraining הרגיל של מודלי שפה הוא שימוש שונה בחיזוי הטוקן הבא - לא דרך סופטמקס אלא תגמול בינארי. המאמר כמובן מוכיח אמפירית שזה משפר את ביצועי המודל. ", 2)

# 7) "two additional models..." paragraph -> short closing remark
$null = $d.Content.Find.Execute("שני מודלים נוספים הם הדקודרים עם משקלים משותפים(בעלי שכבה אחת בלבד כל אחד). הראשון לוקח את הפלט של אנקודר הפלט ומעביר אותו בחזרה למרחב המקורי (עם לוס ריבועי למשל). הדקודר השני לוקח את חיזוי עבור הפלט ומעביר אותם לרחב המקורי של הפלט. ", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר נחמד - קריאה קלילה לסופ״ש….", 2)

# 8) Remove the old closing-recommendation paragraph entirely (its own paragraph, deleted outright)
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "וזה וזה - מאמר נחמד ולא רגיל, מומלץ בחום") {
        $p.Range.Delete()
    }
}

# 9) New arxiv link replaces the old one
$null = $d.Content.Find.Execute("https://arxiv.org/pdf/2502.02996", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2506.08007", 2)
